# Applies the "make forex gains from dividends tax-free" edit to the
# "Foreign Currencies" sheet (rows 3 & 4 relate to FOREX that was never
# acquired - e.g. it came from dividend payments - so it is not a taxable
# disposal: gains are zeroed out and the comment explains why), a couple
# of related value corrections on that same sheet, and the downstream
# "ELSTER - Summary" total that mirrors the Foreign Currencies totals.

$wb = $excel.ActiveWorkbook

$fx = $wb.Worksheets.Item("Foreign Currencies")

# Row 2 - amount corrected slightly
$fx.Range("B2").Value = 1247.91

# Row 3 - FOREX not acquired (received as e.g. dividend payment) -> no taxable gain
$fx.Range("G3").Value = 0
$fx.Range("H3").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 4 - same reasoning as row 3
$fx.Range("G4").Value = 0
$fx.Range("H4").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 5 - amount + gain corrected
$fx.Range("B5").Value = 2567.09
$fx.Range("G5").Value = 20.6

# Row 6 - amount corrected
$fx.Range("B6").Value = 849.87

# Row 7 - amount + gain corrected
$fx.Range("B7").Value = 135.13
$fx.Range("G7").Value = -8.83

# Summary rows recomputed given the above changes
$fx.Range("G9").Value = -43.75
$fx.Range("G10").Value = 20.6
$fx.Range("G11").Value = -64.34999999999999

# The "ELSTER - Summary" sheet duplicates the Foreign Currencies
# "Gains (incl. losses)" total (row 9 / G column above) in its own row 7.
$elster = $wb.Worksheets.Item("ELSTER - Summary")
$elster.Range("C7").Value = -43.75
